$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the cell values can be updated
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure text (cell A10)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-03 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.4840602625630518
$ws.Range("E2").Value = 0.002331002331002363

$ws.Range("D3").Value = 0.3334830370807522
$ws.Range("E3").Value = 0.01051413137841228

$ws.Range("D4").Value = 0.09798543295932143
$ws.Range("E4").Value = -0.002298037829238209

$ws.Range("D5").Value = 0.05409130443927287
$ws.Range("E5").Value = -0.002973467520585471

$ws.Range("D6").Value = 0.03037996295760167
$ws.Range("E6").Value = 0.009163103237629588

$ws.Range("E7").Value = 0.004526991833103899

# Restore sheet protection to its prior (protected) state
$ws.Protect()
